# #5: insurance, claim, debt, investment done
#
# Rework the "保險" (insurance) sheet (sheet8, the 8th worksheet) so it
# follows the same company/name/owner/property_category/category/date/
# legislator_name/legislator_id/source_file/index column layout already
# used by the other property-type sheets (land/building/car/deposit/
# stock/fund/otherbonds). The old free-text "insurance period" column
# (E) is dropped and replaced by the shared metadata columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(8)

# --- Header row ---------------------------------------------------------
$ws.Range("B1").Value = "company"
$ws.Range("C1").Value = "name"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "property_category"
$ws.Range("F1").Value = "category"
$ws.Range("G1").Value = "date"
$ws.Range("H1").Value = "legislator_name"
$ws.Range("I1").Value = "legislator_id"
$ws.Range("J1").Value = "source_file"
$ws.Range("K1").Value = "index"

# --- Data rows -----------------------------------------------------------
# A (index) keeps its original value; B/C/D (company/name/owner) keep
# their original values too - only the tail columns change.
$rows = @(
    @{ Row=2;  A=132; B="南山人壽"; C="南山伴我一生變額壽險";           D="洪佳君" },
    @{ Row=3;  A=133; B="南山人壽"; C="南山好吉利21年期還本養老保險";   D="洪佳君" },
    @{ Row=4;  A=134; B="南山人壽"; C="南山新新增額養老保險";           D="洪佳君" },
    @{ Row=5;  A=135; B="南山人壽"; C="南山好吉利21年期還本養老保險";   D="洪佳君" },
    @{ Row=6;  A=136; B="中泰人壽"; C="中泰人壽金富貴外幣變額年金保險"; D="洪佳君" },
    @{ Row=7;  A=137; B="安聯人壽"; C="世界觀外幣變額萬能壽險";         D="洪佳君" },
    @{ Row=8;  A=138; B="國泰人壽"; C="國泰美滿人生312終身壽險";       D="洪佳君" },
    @{ Row=9;  A=139; B="中華郵政"; C="六年吉利保險";                   D="洪佳君" },
    @{ Row=10; A=141; B="中華郵政"; C="六年吉利保險";                   D="頁志雄" },
    @{ Row=11; A=142; B="中華郵政"; C="六年吉利保險";                   D="黃志雄" },
    @{ Row=12; A=143; B="新光人壽"; C="新光人壽全心终身還本保險";       D="黃志雄" },
    @{ Row=13; A=144; B="新光人壽"; C="新光人壽全意终身還本保險";       D="黃志雄" },
    @{ Row=14; A=145; B="南山人壽"; C="南山人壽鴻利發還本終身分紅保險"; D="貝志雄" }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Range("B$row").Value = $r.B
    $ws.Range("C$row").Value = $r.C
    $ws.Range("D$row").Value = $r.D
    $ws.Range("E$row").Value = "insurance"
    $ws.Range("F$row").Value = "normal"
    # Assigning a plain "2011-12-21" to .Value would get auto-parsed into a
    # date serial number; a leading apostrophe keeps it as the literal text
    # shared-string the source data uses.
    $ws.Range("G$row").Value = "'2011-12-21"
    $ws.Range("H$row").Value = "黃志雄"
    $ws.Cells.Item($row, 9).Value = 1366
    $ws.Range("J$row").Value = "tmp51f51"
    $ws.Cells.Item($row, 11).Value = $r.A
}
